$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_5a_Indikatoren")

# Row 7
$ws.Range("E7").Value = "XVermeidbare Sterblichkeit"
$ws.Range("G7").Value = "XDurch Prävention und Behandlung vermeidbare Sterblichkeit"
$ws.Range("I7").Value = "XSenkung auf 200 Todesfälle je 100 000 Einwohnerinnen und Einwohner bis 2030"
$ws.Range("K7").Value = "XSenkung auf 200 Todesfälle je 100 000 Einwohnerinnen und Einwohner bis 2030"

# Row 15
$ws.Range("E15").Value = "XXXFortschritte beim Global Health Security Index"
$ws.Range("G15").Value = "XAnteil der Länder, die Deutschland bei der XXXPandemieprävention, -vorsorge und -reaktion (PPR) unterstützt und die im Global Health Security Index (GHSI) Fortschritte gemacht haben, an der Anzahl der Länder, die Deutschland im Bereich PPR unterstützt"
$ws.Range("I15").Value = "XXX30 Prozent der Länder haben im Jahr 2030 eine Verbesserung des Global Health"
$ws.Range("J15").Value = "30 Prozent der Länder haben im Jahr 2030 eine Verbesserung des Global Health Security Index (GSHI) um mindestens den Wert 3 gegenüber 2019 erzielt"
$ws.Range("K15").Value = "XXX30 Prozent der Länder haben im Jahr 2030 eine Verbesserung des Global Health Security Index (GSHI) um mindestens den Wert 3 gegenüber 2019 erzielt"

# Row 16
$ws.Range("D16").Value = "Unterschied in der Lebenserwartung zwischen sozio-ökonomisch deprivierten und wohlhabenden Regionen (Frauen)"
$ws.Range("E16").Value = "XXXUnterschied in der Lebenserwartung zwischen sozio-ökonomisch deprivierten und wohlhabenden Regionen (Frauen)"
$ws.Range("G16").Value = "XXXUnterschied in der Lebenserwartung zwischen sozio-ökonomisch deprivierten und wohlhabenden Regionen (Frauen)"
$ws.Range("I16").Value = "XXXRückgang des Unterschieds in der mittleren Lebenserwartung von Frauen zwischen deprivierten und wohlhabenden Regionen bei gleichzeitigem Anstieg der Lebenserwartung von Frauen in deprivierten Regionen"
$ws.Range("J16").Value = "ückgang des Unterschieds in der mittleren Lebenserwartung von Frauen zwischen deprivierten und wohlhabenden Regionen bei gleichzeitigem Anstieg der Lebenserwartung von Frauen in deprivierten Regionen"
$ws.Range("K16").Value = "XXXückgang des Unterschieds in der mittleren Lebenserwartung von Frauen zwischen deprivierten und wohlhabenden Regionen bei gleichzeitigem Anstieg der Lebenserwartung von Frauen in deprivierten Regionen"

# Row 17
$ws.Range("D17").Value = "Unterschied in der Lebenserwartung zwischen sozio-ökonomisch deprivierten und wohlhabenden Regionen (Männer)"
$ws.Range("E17").Value = "XXXUnterschied in der Lebenserwartung zwischen sozio-ökonomisch deprivierten und wohlhabenden Regionen (Männer)"
$ws.Range("G17").Value = "XXXUnterschied in der Lebenserwartung zwischen sozio-ökonomisch deprivierten und wohlhabenden Regionen (Männer)"
$ws.Range("I17").Value = "XXXRückgang des Unterschieds in der mittleren Lebenserwartung von Männern zwischen deprivierten und wohlhabenden Regionen bei gleichzeitigem Anstieg der Lebenserwartung von Frauen in deprivierten Regionen"
$ws.Range("J17").Value = "Rückgang des Unterschieds in der mittleren Lebenserwartung von Männern zwischen deprivierten und wohlhabenden Regionen bei gleichzeitigem Anstieg der Lebenserwartung von Frauen in deprivierten Regionen"
$ws.Range("K17").Value = "XXXRückgang des Unterschieds in der mittleren Lebenserwartung von Männern zwischen deprivierten und wohlhabenden Regionen bei gleichzeitigem Anstieg der Lebenserwartung von Frauen in deprivierten Regionen"

# Row 20
$ws.Range("I20").Value = "XXXKontinuierlicher Anstieg"
